# Add the new employee record (row 5) to the Working Hours sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# employee_id "69" must be stored as text (it is not a numeric id in the
# source data), so force the text number format before assigning the
# value, then clear the format again so the cell keeps the default style
# (matches how the rest of the data rows are styled) while remaining text.
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "69"
$ws.Range("A5").ClearFormats()

$ws.Range("B5").Value = "Thim"
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = 24
$ws.Range("E5").Value = 0
